$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in grades for students that previously had "Belum Mengerjakan Tugas Akhir"
$ws.Range("D2").Value = 90
$ws.Range("D15").Value = 82
$ws.Range("D24").Value = 92

# Update an existing grade
$ws.Range("D28").Value = 94

# Update the view: scroll back to top and change selection
$ws.Activate()
$ws.Range("D11").Select()
$excel.ActiveWindow.ScrollRow = 1
